# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" between "2021-Q2" and "总计", holding the
#    per-fund holding detail for the new quarter.
# 2. Update the "总计" (summary) sheet with a new top row for 2022-Q1,
#    pushing the existing 2021-Q2 summary row down.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item("2021-Q2")
$totalSheet = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------------
# Step 1: create "2022-Q1" by duplicating "总计" (so it inherits the exact same
# header / index-column formatting) and placing it right after "2021-Q2".
# ---------------------------------------------------------------------------
$totalSheet.Copy($null, $sheet1)
$newSheet = $wb.Worksheets.Item("总计 (2)")
$newSheet.Name = "2022-Q1"

# The duplicated sheet only carries style for columns B:D (and index col A for
# row 2) - extend the same formatting to the extra columns (E:H) and extra
# index rows (3:4) this sheet needs, by copying formats within the sheet.
$newSheet.Range("B1").Copy()
$newSheet.Range("E1:H1").PasteSpecial(-4122)
$newSheet.Range("A2").Copy()
$newSheet.Range("A3:A4").PasteSpecial(-4122)

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Fund holding rows. The numeric-looking columns (code / size / position /
# ratio / value) are stored as text in the source file, so they are entered
# with a leading apostrophe to force text, matching the original layout.
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'001877"
$newSheet.Range("C2").Value = "宝盈国家安全战略沪港深股票"
$newSheet.Range("D2").Value = "'13.59"
$newSheet.Range("E2").Value = "'90.66"
$newSheet.Range("F2").Value = "'4.41"
$newSheet.Range("G2").Value = "'0.5993"
$newSheet.Range("H2").Value = 4

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'001487"
$newSheet.Range("C3").Value = "宝盈优势产业灵活配置混合"
$newSheet.Range("D3").Value = "'17.02"
$newSheet.Range("E3").Value = "'91.61"
$newSheet.Range("F3").Value = "'3.14"
$newSheet.Range("G3").Value = "'0.5344"
$newSheet.Range("H3").Value = 6

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "'001075"
$newSheet.Range("C4").Value = "宝盈转型动力灵活配置混合"
$newSheet.Range("D4").Value = "'5.13"
$newSheet.Range("E4").Value = "'86.64"
$newSheet.Range("F4").Value = "'3.32"
$newSheet.Range("G4").Value = "'0.1703"
$newSheet.Range("H4").Value = 7

# The leading-apostrophe entry leaves a "quote prefix" style marker behind;
# strip it back to the plain (no explicit style) look of the source file by
# pasting in the format of an untouched default cell, then clean that helper
# cell up again.
$newSheet.Range("Z99").Value = "x"
$newSheet.Range("Z99").Copy()
$newSheet.Range("B2:G2").PasteSpecial(-4122)
$newSheet.Range("B3:G3").PasteSpecial(-4122)
$newSheet.Range("B4:G4").PasteSpecial(-4122)
$newSheet.Range("Z99").ClearContents()

# ---------------------------------------------------------------------------
# Step 2: update "总计" - push the existing 2021-Q2 summary row to row 3 and
# write the new 2022-Q1 totals into row 2.
# Re-fetch the sheet by name: inserting/copying sheets above shifts indices,
# and the old $totalSheet handle would otherwise now resolve to "2022-Q1".
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Range("A2").Copy()
$totalSheet.Range("A3").PasteSpecial(-4122)

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q2"
$totalSheet.Range("C3").Value = 2
$totalSheet.Range("D3").Value = 0.07000000000000001

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 3
$totalSheet.Range("D2").Value = 1.3

Write-Host "Done: 2022-Q1 sheet added and summary sheet updated"
